$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Add new "on_topic" column (column I) -------------------------------
# Header
$ws.Cells.Item(1, 9).Value = "on_topic"

# Rows that represent an actual trial (i.e. already have a value in column F)
# get "Y" in the new on_topic column.
$trialRows = @(38,44,49,52,59,63,69,73,78,82,83,90,92,97,99,104,105,110,114,115,116,120,124,127,131,133)
foreach ($r in $trialRows) {
    $ws.Cells.Item($r, 9).Value = "Y"
}

# --- View / window state --------------------------------------------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.Zoom = 150

# Re-establish the frozen pane (row 1) and scroll the frozen region so the
# same relative window (row 117 downward) is showing, then restore the
# previously-active selection at I134.
$win.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$win.FreezePanes = $true
$win.ScrollRow = 117
$win.ScrollColumn = 1

$ws.Range("A1").Select() | Out-Null
$ws.Range("I134").Select() | Out-Null
